# Auto-generated PowerShell Excel COM-interop script
# Applies cell value updates across multiple worksheets per the target diff
# Commit message: Add data for 2024-10-18

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals (18 cell updates)
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6136
$ws.Range('K3').Value = 6314
$ws.Range('K4').Value = 368
$ws.Range('K5').Value = 5846
$ws.Range('D6').Value = 1614
$ws.Range('F6').Value = 1580
$ws.Range('K6').Value = 1325
$ws.Range('K7').Value = 450
$ws.Range('K8').Value = 16902
$ws.Range('K9').Value = 6960
$ws.Range('I10').Value = 43536
$ws.Range('J10').Value = 45415
$ws.Range('K10').Value = 45675
$ws.Range('D11').Value = 93462
$ws.Range('F11').Value = 84551
$ws.Range('I11').Value = 84926
$ws.Range('J11').Value = 98625
$ws.Range('K11').Value = 89976

# Sheet 2: By Neighborhood (95 cell updates)
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 856
$ws.Range('K4').Value = 472
$ws.Range('K6').Value = 846
$ws.Range('K7').Value = 1935
$ws.Range('K8').Value = 3743
$ws.Range('K9').Value = 358
$ws.Range('K10').Value = 988
$ws.Range('K11').Value = 1642
$ws.Range('K12').Value = 436
$ws.Range('K14').Value = 532
$ws.Range('K15').Value = 739
$ws.Range('K16').Value = 702
$ws.Range('K17').Value = 107
$ws.Range('K18').Value = 539
$ws.Range('K19').Value = 1967
$ws.Range('K20').Value = 1500
$ws.Range('K22').Value = 350
$ws.Range('K23').Value = 1007
$ws.Range('K24').Value = 494
$ws.Range('K25').Value = 401
$ws.Range('K27').Value = 1167
$ws.Range('K29').Value = 2947
$ws.Range('K31').Value = 789
$ws.Range('K32').Value = 163
$ws.Range('K33').Value = 2171
$ws.Range('K34').Value = 752
$ws.Range('K35').Value = 267
$ws.Range('K36').Value = 1269
$ws.Range('K37').Value = 2064
$ws.Range('K38').Value = 131
$ws.Range('K39').Value = 150
$ws.Range('K40').Value = 267
$ws.Range('K41').Value = 374
$ws.Range('K42').Value = 2425
$ws.Range('K43').Value = 958
$ws.Range('K44').Value = 999
$ws.Range('K46').Value = 326
$ws.Range('K47').Value = 771
$ws.Range('K48').Value = 2503
$ws.Range('K49').Value = 1634
$ws.Range('K50').Value = 789
$ws.Range('K51').Value = 1114
$ws.Range('K52').Value = 1448
$ws.Range('K53').Value = 1735
$ws.Range('K54').Value = 3197
$ws.Range('K55').Value = 993
$ws.Range('K56').Value = 543
$ws.Range('K57').Value = 427
$ws.Range('K59').Value = 199
$ws.Range('K60').Value = 591
$ws.Range('K61').Value = 126
$ws.Range('D63').Value = 1066
$ws.Range('F63').Value = 1044
$ws.Range('I63').Value = 1576
$ws.Range('J63').Value = 477
$ws.Range('K63').Value = 380
$ws.Range('K64').Value = 820
$ws.Range('K65').Value = 1281
$ws.Range('K66').Value = 552
$ws.Range('K67').Value = 1859
$ws.Range('K68').Value = 366
$ws.Range('K70').Value = 519
$ws.Range('K71').Value = 283
$ws.Range('K72').Value = 573
$ws.Range('K73').Value = 1153
$ws.Range('K74').Value = 279
$ws.Range('K75').Value = 350
$ws.Range('K76').Value = 2249
$ws.Range('K77').Value = 344
$ws.Range('K78').Value = 1552
$ws.Range('K79').Value = 1798
$ws.Range('K80').Value = 372
$ws.Range('K83').Value = 1270
$ws.Range('K84').Value = 656
$ws.Range('K85').Value = 2959
$ws.Range('K86').Value = 800
$ws.Range('K87').Value = 365
$ws.Range('K88').Value = 760
$ws.Range('K89').Value = 1803
$ws.Range('K90').Value = 872
$ws.Range('K91').Value = 782
$ws.Range('K92').Value = 284
$ws.Range('K93').Value = 703
$ws.Range('K94').Value = 2306
$ws.Range('K95').Value = 1088
$ws.Range('K96').Value = 1255
$ws.Range('K97').Value = 1353
$ws.Range('K98').Value = 1144
$ws.Range('K99').Value = 1169
$ws.Range('K100').Value = 260
$ws.Range('D101').Value = 93462
$ws.Range('F101').Value = 84551
$ws.Range('I101').Value = 84926
$ws.Range('J101').Value = 98625
$ws.Range('K101').Value = 89976

# Sheet 3: Bridgeport (2 cell updates)
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K10').Value = 281
$ws.Range('K11').Value = 532

# Sheet 4: West Ridge (5 cell updates)
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 69
$ws.Range('K8').Value = 322
$ws.Range('K9').Value = 95
$ws.Range('K10').Value = 622
$ws.Range('K11').Value = 1255

# Sheet 5: Auburn Gresham (3 cell updates)
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K3').Value = 205
$ws.Range('K10').Value = 674
$ws.Range('K11').Value = 1935

# Sheet 6: Belmont Cragin (3 cell updates)
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K5').Value = 80
$ws.Range('K10').Value = 856
$ws.Range('K11').Value = 1642

# Sheet 7: O'Hare (3 cell updates)
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('K8').Value = 95
$ws.Range('K10').Value = 379
$ws.Range('K11').Value = 519

# Sheet 8: Uptown (5 cell updates)
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 89
$ws.Range('K3').Value = 95
$ws.Range('K6').Value = 38
$ws.Range('K10').Value = 1131
$ws.Range('K11').Value = 1803

# Sheet 9: South Shore (7 cell updates)
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 321
$ws.Range('K3').Value = 339
$ws.Range('K5').Value = 352
$ws.Range('K8').Value = 500
$ws.Range('K9').Value = 239
$ws.Range('K10').Value = 1107
$ws.Range('K11').Value = 2959

# Sheet 10: Little Village (3 cell updates)
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K3').Value = 161
$ws.Range('K10').Value = 661
$ws.Range('K11').Value = 1448

# Sheet 12: Logan Square (4 cell updates)
$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K5').Value = 96
$ws.Range('K8').Value = 272
$ws.Range('K10').Value = 1091
$ws.Range('K11').Value = 1735

# Sheet 13: Austin (4 cell updates)
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 425
$ws.Range('K8').Value = 742
$ws.Range('K10').Value = 1419
$ws.Range('K11').Value = 3743

# Sheet 14: Jefferson Park (2 cell updates)
$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('K10').Value = 185
$ws.Range('K11').Value = 326

# Sheet 15: Morgan Park (3 cell updates)
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K8').Value = 146
$ws.Range('K10').Value = 286
$ws.Range('K11').Value = 591

# Sheet 16: Oakland (3 cell updates)
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K3').Value = 21
$ws.Range('K10').Value = 114
$ws.Range('K11').Value = 283

# Sheet 17: South Chicago (3 cell updates)
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 160
$ws.Range('K8').Value = 258
$ws.Range('K11').Value = 1270

# Sheet 18: Garfield Park (7 cell updates)
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 239
$ws.Range('K3').Value = 332
$ws.Range('K6').Value = 47
$ws.Range('K8').Value = 408
$ws.Range('K9').Value = 282
$ws.Range('K10').Value = 704
$ws.Range('K11').Value = 2171

# Sheet 19: Roseland (6 cell updates)
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 176
$ws.Range('K5').Value = 133
$ws.Range('K8').Value = 441
$ws.Range('K9').Value = 134
$ws.Range('K10').Value = 683
$ws.Range('K11').Value = 1798

# Sheet 20: Pullman (2 cell updates)
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K10').Value = 204
$ws.Range('K11').Value = 350

# Sheet 21: Hegewisch (2 cell updates)
$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('K10').Value = 131
$ws.Range('K11').Value = 267

# Sheet 22: Near South Side (3 cell updates)
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K9').Value = 50
$ws.Range('K10').Value = 459
$ws.Range('K11').Value = 820

# Sheet 23: West Pullman (2 cell updates)
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K10').Value = 402
$ws.Range('K11').Value = 1088

# Sheet 24: Grand Crossing (4 cell updates)
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 239
$ws.Range('K8').Value = 458
$ws.Range('K10').Value = 696
$ws.Range('K11').Value = 2064

# Sheet 25: Edgewater (5 cell updates)
$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K3').Value = 46
$ws.Range('K5').Value = 74
$ws.Range('K6').Value = 27
$ws.Range('K10').Value = 713
$ws.Range('K11').Value = 1167

# Sheet 26: New City (4 cell updates)
$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 168
$ws.Range('K9').Value = 178
$ws.Range('K10').Value = 470
$ws.Range('K11').Value = 1281

# Sheet 27: Woodlawn (4 cell updates)
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K4').Value = 6
$ws.Range('K8').Value = 235
$ws.Range('K10').Value = 483
$ws.Range('K11').Value = 1169

# Sheet 29: Gage Park (2 cell updates)
$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K8').Value = 168
$ws.Range('K11').Value = 789

# Sheet 30: North Lawndale (7 cell updates)
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 230
$ws.Range('K3').Value = 300
$ws.Range('K7').Value = 20
$ws.Range('K8').Value = 309
$ws.Range('K9').Value = 238
$ws.Range('K10').Value = 608
$ws.Range('K11').Value = 1859

# Sheet 31: South Deering (3 cell updates)
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K9').Value = 32
$ws.Range('K10').Value = 310
$ws.Range('K11').Value = 656

# Sheet 32: West Loop (3 cell updates)
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K9').Value = 126
$ws.Range('K10').Value = 1615
$ws.Range('K11').Value = 2306

# Sheet 33: River North (3 cell updates)
$ws = $wb.Worksheets.Item('River North')
$ws.Range('K9').Value = 148
$ws.Range('K10').Value = 1639
$ws.Range('K11').Value = 2249

# Sheet 34: Ukrainian Village (3 cell updates)
$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('K5').Value = 25
$ws.Range('K10').Value = 227
$ws.Range('K11').Value = 365

# Sheet 35: East Side (2 cell updates)
$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K10').Value = 143
$ws.Range('K11').Value = 401

# Sheet 36: Wrigleyville (2 cell updates)
$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('K3').Value = 6
$ws.Range('K11').Value = 260

# Sheet 37: Bucktown (3 cell updates)
$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K5').Value = 63
$ws.Range('K10').Value = 486
$ws.Range('K11').Value = 702

# Sheet 38: Lincoln Park (4 cell updates)
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K8').Value = 168
$ws.Range('K9').Value = 58
$ws.Range('K10').Value = 1222
$ws.Range('K11').Value = 1634

# Sheet 39: West Town (4 cell updates)
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K8').Value = 333
$ws.Range('K9').Value = 92
$ws.Range('K10').Value = 730
$ws.Range('K11').Value = 1353

# Sheet 40: Lower West Side (4 cell updates)
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K4').Value = 5
$ws.Range('K8').Value = 219
$ws.Range('K10').Value = 472
$ws.Range('K11').Value = 993

# Sheet 41: Loop (4 cell updates)
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K5').Value = 65
$ws.Range('K6').Value = 20
$ws.Range('K10').Value = 2347
$ws.Range('K11').Value = 3197

# Sheet 42: Portage Park (2 cell updates)
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K10').Value = 662
$ws.Range('K11').Value = 1153

# Sheet 43: Englewood (7 cell updates)
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 331
$ws.Range('K3').Value = 413
$ws.Range('K5').Value = 148
$ws.Range('K6').Value = 56
$ws.Range('K9').Value = 322
$ws.Range('K10').Value = 956
$ws.Range('K11').Value = 2947

# Sheet 44: Lake View (5 cell updates)
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 38
$ws.Range('K5').Value = 187
$ws.Range('K6').Value = 38
$ws.Range('K10').Value = 1739
$ws.Range('K11').Value = 2503

# Sheet 45: Chatham (4 cell updates)
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 182
$ws.Range('K3').Value = 187
$ws.Range('K10').Value = 757
$ws.Range('K11').Value = 1967

# Sheet 46: North Center (2 cell updates)
$ws = $wb.Worksheets.Item('North Center')
$ws.Range('K10').Value = 305
$ws.Range('K11').Value = 552

# Sheet 47: Irving Park (2 cell updates)
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K9').Value = 71
$ws.Range('K11').Value = 999

# Sheet 48: Humboldt Park (6 cell updates)
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 216
$ws.Range('K5').Value = 186
$ws.Range('K8').Value = 523
$ws.Range('K9').Value = 293
$ws.Range('K10').Value = 912
$ws.Range('K11').Value = 2425

# Sheet 49: Clearing (2 cell updates)
$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('K8').Value = 106
$ws.Range('K11').Value = 350

# Sheet 50: Ashburn (3 cell updates)
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 55
$ws.Range('K10').Value = 315
$ws.Range('K11').Value = 846

# Sheet 51: Hermosa (2 cell updates)
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K10').Value = 122
$ws.Range('K11').Value = 374

# Sheet 52: Grand Boulevard (4 cell updates)
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K2').Value = 107
$ws.Range('K8').Value = 354
$ws.Range('K10').Value = 554
$ws.Range('K11').Value = 1269

# Sheet 54: Avondale (2 cell updates)
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K10').Value = 677
$ws.Range('K11').Value = 988

# Sheet 55: Streeterville (3 cell updates)
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K9').Value = 32
$ws.Range('K10').Value = 568
$ws.Range('K11').Value = 800

# Sheet 56: Rogers Park (3 cell updates)
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K8').Value = 284
$ws.Range('K10').Value = 936
$ws.Range('K11').Value = 1552

# Sheet 57: North Park (3 cell updates)
$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K8').Value = 55
$ws.Range('K10').Value = 226
$ws.Range('K11').Value = 366

# Sheet 58: Dunning (2 cell updates)
$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K10').Value = 273
$ws.Range('K11').Value = 494

# Sheet 59: Brighton Park (2 cell updates)
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K10').Value = 318
$ws.Range('K11').Value = 739

# Sheet 60: Douglas (2 cell updates)
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K10').Value = 543
$ws.Range('K11').Value = 1007

# Sheet 61: Washington Park (3 cell updates)
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K8').Value = 227
$ws.Range('K10').Value = 273
$ws.Range('K11').Value = 782

# Sheet 62: Little Italy, UIC (4 cell updates)
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 73
$ws.Range('K9').Value = 88
$ws.Range('K10').Value = 566
$ws.Range('K11').Value = 1114

# Sheet 63: Chicago Lawn (4 cell updates)
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 171
$ws.Range('K5').Value = 101
$ws.Range('K10').Value = 556
$ws.Range('K11').Value = 1500

# Sheet 65: Kenwood (3 cell updates)
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K8').Value = 222
$ws.Range('K10').Value = 366
$ws.Range('K11').Value = 771

# Sheet 66: Washington Heights (2 cell updates)
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 73
$ws.Range('K11').Value = 872

# Sheet 67: Lincoln Square (3 cell updates)
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K8').Value = 141
$ws.Range('K10').Value = 472
$ws.Range('K11').Value = 789

# Sheet 68: West Lawn (3 cell updates)
$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('K8').Value = 224
$ws.Range('K10').Value = 350
$ws.Range('K11').Value = 703

# Sheet 69: Calumet Heights (2 cell updates)
$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K10').Value = 215
$ws.Range('K11').Value = 539

# Sheet 70: Riverdale (2 cell updates)
$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K8').Value = 72
$ws.Range('K11').Value = 344

# Sheet 71: Magnificent Mile (2 cell updates)
$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range('K10').Value = 509
$ws.Range('K11').Value = 543

# Sheet 73: Albany Park (3 cell updates)
$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K3').Value = 48
$ws.Range('K10').Value = 476
$ws.Range('K11').Value = 856

# Sheet 74: Old Town (3 cell updates)
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K6').Value = 8
$ws.Range('K10').Value = 367
$ws.Range('K11').Value = 573

# Sheet 75: Hyde Park (2 cell updates)
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K10').Value = 552
$ws.Range('K11').Value = 958

# Sheet 76: Burnside (2 cell updates)
$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('K10').Value = 25
$ws.Range('K11').Value = 107

# Sheet 77: Archer Heights (2 cell updates)
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('K10').Value = 258
$ws.Range('K11').Value = 472

# Sheet 78: Rush & Division (2 cell updates)
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K5').Value = 18
$ws.Range('K11').Value = 372

# Sheet 79: Garfield Ridge (2 cell updates)
$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K10').Value = 303
$ws.Range('K11').Value = 752

# Sheet 80: Gold Coast (2 cell updates)
$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range('K9').Value = 205
$ws.Range('K10').Value = 267

# Sheet 82: Galewood (3 cell updates)
$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('K8').Value = 41
$ws.Range('K10').Value = 89
$ws.Range('K11').Value = 163

# Sheet 84: Wicker Park (3 cell updates)
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K8').Value = 191
$ws.Range('K10').Value = 760
$ws.Range('K11').Value = 1144

# Sheet 86: West Elsdon (2 cell updates)
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('K10').Value = 105
$ws.Range('K11').Value = 284

# Sheet 87: Mount Greenwood (2 cell updates)
$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('K9').Value = 84
$ws.Range('K10').Value = 126

# Sheet 89: Avalon Park (3 cell updates)
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K5').Value = 27
$ws.Range('K10').Value = 140
$ws.Range('K11').Value = 358

# Sheet 90: Montclare (2 cell updates)
$ws = $wb.Worksheets.Item('Montclare')
$ws.Range('K10').Value = 101
$ws.Range('K11').Value = 199

# Sheet 91: Greektown (2 cell updates)
$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('K4').Value = 13
$ws.Range('K9').Value = 150

# Sheet 92: Mckinley Park (2 cell updates)
$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K10').Value = 254
$ws.Range('K11').Value = 427

# Sheet 93: United Center (3 cell updates)
$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K9').Value = 94
$ws.Range('K10').Value = 344
$ws.Range('K11').Value = 760

# Sheet 94: Printers Row (2 cell updates)
$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range('K10').Value = 230
$ws.Range('K11').Value = 279

# Sheet 97: Grant Park (2 cell updates)
$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('K9').Value = 104
$ws.Range('K10').Value = 131

# Sheet 98: Beverly (2 cell updates)
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('K10').Value = 321
$ws.Range('K11').Value = 436
